$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 9473.666999999999

$ws.Range("H62").Value = 3913.8572
$ws.Range("I62").Value = 3647.75
$ws.Range("K62").Value = 3647.75
$ws.Range("M62").Value = -3023.75

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H65").Value = 3913.8572
$ws.Range("I65").Value = 3647.75
$ws.Range("K65").Value = 18238.75
$ws.Range("M65").Value = -15118.75

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H138").Value = 2534.0857
$ws.Range("J138").Value = 2206.8708
$ws.Range("L138").Value = 6620.6124
$ws.Range("N138").Value = -16900.6124

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2033.1666
$ws.Range("I2").Value = 2139.8
$ws.Range("K2").Value = 2139.8
$ws.Range("M2").Value = -2026.8

$ws.Range("H61").Value = 3927.8
$ws.Range("I61").Value = 3963.5
$ws.Range("J61").Value = 3874.25
$ws.Range("K61").Value = 3963.5
$ws.Range("L61").Value = 3874.25
$ws.Range("M61").Value = -3751.5
$ws.Range("N61").Value = -4298.25

$ws.Range("H63").Value = 6212
$ws.Range("I63").Value = 5328
$ws.Range("J63").Value = 7096
$ws.Range("K63").Value = 5328
$ws.Range("L63").Value = 7096
$ws.Range("M63").Value = -4642
$ws.Range("N63").Value = -8468

$ws.Range("H66").Value = 6212
$ws.Range("I66").Value = 5328
$ws.Range("J66").Value = 7096
$ws.Range("K66").Value = 26640
$ws.Range("L66").Value = 35480
$ws.Range("M66").Value = -23208
$ws.Range("N66").Value = -42344

$ws.Range("H74").Value = 2392.625
$ws.Range("I74").Value = 2392.625
$ws.Range("K74").Value = 2392.625
$ws.Range("M74").Value = -1518.625

$ws.Range("H77").Value = 2392.625
$ws.Range("I77").Value = 2392.625
$ws.Range("K77").Value = 11963.125
$ws.Range("M77").Value = -7595.125

$ws.Range("H110").Value = 2057.25
$ws.Range("J110").Value = 1909
$ws.Range("L110").Value = 1909
$ws.Range("N110").Value = -5999

$ws.Range("H116").Value = 2033.1666
$ws.Range("I116").Value = 2139.8
$ws.Range("K116").Value = 2139.8
$ws.Range("M116").Value = 154.1999999999998

$ws.Range("H136").Value = 3927.8
$ws.Range("I136").Value = 3963.5
$ws.Range("J136").Value = 3874.25
$ws.Range("K136").Value = 11890.5
$ws.Range("L136").Value = 11622.75
$ws.Range("M136").Value = -9340.5
$ws.Range("N136").Value = -16722.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2033.1666
$ws.Range("I3").Value = 2139.8
$ws.Range("K3").Value = 2139.8
$ws.Range("M3").Value = -2025.8

$ws.Range("H22").Value = 466.83334
$ws.Range("I22").Value = 359.8
$ws.Range("K22").Value = 359.8
$ws.Range("M22").Value = -186.8

$ws.Range("H105").Value = 3998.5715
$ws.Range("I105").Value = 3668.6667
$ws.Range("J105").Value = 4246
$ws.Range("K105").Value = 3668.6667
$ws.Range("L105").Value = 4246
$ws.Range("M105").Value = -1921.6667
$ws.Range("N105").Value = -7740

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 468.2
$ws.Range("J22").Value = 419.5
$ws.Range("L22").Value = 419.5
$ws.Range("N22").Value = -1119.5

$ws.Range("H58").Value = 3432.4443
$ws.Range("J58").Value = 3897
$ws.Range("L58").Value = 3897
$ws.Range("N58").Value = -4303

$ws.Range("H136").Value = 3432.4443
$ws.Range("J136").Value = 3897
$ws.Range("L136").Value = 11691
$ws.Range("N136").Value = -16791

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 1578.4
$ws.Range("I141").Value = 1578.4
$ws.Range("K141").Value = 4735.200000000001
$ws.Range("M141").Value = 444.7999999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9451.388999999999
$ws.Range("I80").Value = 4517.6
$ws.Range("J80").Value = 15618.625
$ws.Range("K80").Value = 4517.6
$ws.Range("L80").Value = 15618.625
$ws.Range("M80").Value = -3519.6
$ws.Range("N80").Value = -17614.625

$ws.Range("H83").Value = 9451.388999999999
$ws.Range("I83").Value = 4517.6
$ws.Range("J83").Value = 15618.625
$ws.Range("K83").Value = 22588
$ws.Range("L83").Value = 78093.125
$ws.Range("M83").Value = -17596
$ws.Range("N83").Value = -88077.125

$ws.Range("H102").Value = 3757.875
$ws.Range("I102").Value = 3161
$ws.Range("K102").Value = 3161
$ws.Range("M102").Value = -1539

$ws.Range("H122").Value = 2760.8
$ws.Range("I122").Value = 1932
$ws.Range("K122").Value = 5796
$ws.Range("M122").Value = -3346

$ws.Range("H136").Value = 43151.625
$ws.Range("J136").Value = 43151.625
$ws.Range("L136").Value = 129454.875
$ws.Range("N136").Value = -134554.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2154.6191
$ws.Range("I7").Value = 1042.4286
$ws.Range("J7").Value = 2710.7144
$ws.Range("K7").Value = 1042.4286
$ws.Range("L7").Value = 2710.7144
$ws.Range("M7").Value = -930.4286
$ws.Range("N7").Value = -2934.7144

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H80").Value = 72000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 72000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws.Range("H126").Value = 2154.6191
$ws.Range("I126").Value = 1042.4286
$ws.Range("J126").Value = 2710.7144
$ws.Range("K126").Value = 3127.2858
$ws.Range("L126").Value = 8132.1432
$ws.Range("M126").Value = -657.2857999999997
$ws.Range("N126").Value = -13072.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5217.636
$ws.Range("J62").Value = 7249.75
$ws.Range("L62").Value = 7249.75
$ws.Range("N62").Value = -8497.75

$ws.Range("H65").Value = 5217.636
$ws.Range("J65").Value = 7249.75
$ws.Range("L65").Value = 36248.75
$ws.Range("N65").Value = -42488.75
